# Updated legacy GSC export data:
# The coverage-over-time export rolled forward by one day - the oldest
# date row (2025-08-31, the first data row under the header) drops off
# and every subsequent day's row shifts up to take its place, with the
# newest day (2025-11-18) now the last row.
#
# Deleting the row (rather than rewriting every cell) reproduces exactly
# that "shift up" semantics, including Excel's usual housekeeping -
# shrinking the sheet dimension by one row and dropping the now-unused
# shared string for the removed date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 is the first data row (2025-08-31). Deleting it shifts every
# later row (2025-09-01 .. 2025-11-18) up by one.
$ws.Rows.Item(2).Delete()
